# save data done + era data updated
# Add a new "Save" column (H) to the sheet with 0/1 values for rows 2-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled the same as the other header cells (B1:G1 use style index 1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Row-by-row "Save" flag values (column H), matching the diff.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 0
    32 = 0
}

foreach ($row in 2..32) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
